# Update crypto price/volume figures to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    # Force Excel to store the value as literal text instead of
    # silently re-interpreting it as a number (which would strip
    # meaningful trailing zeros, introduce floating point noise,
    # or switch to scientific notation for small values).
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "27.212.71"
$ws.Range("E2").Value = "  +0.27%  "
Set-TextValue $ws.Range("D3") "1.895.55"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue $ws.Range("D5") "306.20"
$ws.Range("E5").Value = "  -0.39%  "
$ws.Range("E6").Value = "  -0.06%  "
Set-TextValue $ws.Range("D7") "0.5381"
$ws.Range("E7").Value = "  +3.49%  "
Set-TextValue $ws.Range("D8") "0.3792"
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("E9").Value = "  -0.11%  "
Set-TextValue $ws.Range("D10") "21.99"
$ws.Range("E10").Value = "  +3.81%  "
$ws.Range("E11").Value = "  -0.34%  "
Set-TextValue $ws.Range("D12") "0.08176"
$ws.Range("E12").Value = "  -0.49%  "
Set-TextValue $ws.Range("D13") "94.89"
$ws.Range("E13").Value = "  -1.37%  "
Set-TextValue $ws.Range("D14") "5.343"
$ws.Range("E14").Value = "  +0.34%  "
Set-TextValue $ws.Range("D15") "1.768.29"
$ws.Range("E15").Value = "  -7.07%  "
Set-TextValue $ws.Range("D16") "1.002"
$ws.Range("E16").Value = "  -0.08%  "
Set-TextValue $ws.Range("D17") "14.84"
$ws.Range("E17").Value = "  +1.49%  "
Set-TextValue $ws.Range("D18") "0.000008643"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("E19").Value = "  -0.04%  "
Set-TextValue $ws.Range("D20") "27.074.46"
$ws.Range("E20").Value = "  -0.36%  "
Set-TextValue $ws.Range("D21") "5.034"
$ws.Range("E21").Value = "  -1.18%  "
Set-TextValue $ws.Range("D22") "10.80"
$ws.Range("E22").Value = "  +0.71%  "
Set-TextValue $ws.Range("D23") "6.468"
$ws.Range("E23").Value = "  +0.68%  "
Set-TextValue $ws.Range("D24") "148.67"
$ws.Range("E24").Value = "  -0.20%  "
Set-TextValue $ws.Range("D25") "2.290"
$ws.Range("E25").Value = "  -0.92%  "
Set-TextValue $ws.Range("D26") "18.34"
$ws.Range("E26").Value = "  +0.82%  "
Set-TextValue $ws.Range("D27") "1.758"
$ws.Range("E27").Value = "  +0.83%  "
Set-TextValue $ws.Range("D28") "116.16"
$ws.Range("E28").Value = "  +0.65%  "
Set-TextValue $ws.Range("D29") "4.818"
$ws.Range("E29").Value = "  +0.22%  "
Set-TextValue $ws.Range("D30") "4.645"
$ws.Range("E30").Value = "  -4.38%  "
Set-TextValue $ws.Range("D31") "0.09162"
$ws.Range("E31").Value = "  -0.51%  "
Set-TextValue $ws.Range("D32") "0.8190"
$ws.Range("E32").Value = "  +2.89%  "
$ws.Range("E33").Value = "  +0.22%  "
Set-TextValue $ws.Range("D34") "1.218"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("E35").Value = "  +2.10%  "
Set-TextValue $ws.Range("D36") "3.293"
$ws.Range("E36").Value = "  -3.95%  "
Set-TextValue $ws.Range("D37") "2.672"
$ws.Range("E37").Value = "  +2.51%  "
Set-TextValue $ws.Range("D38") "0.5957"
$ws.Range("E38").Value = "  +4.11%  "
Set-TextValue $ws.Range("D39") "0.01984"
$ws.Range("E39").Value = "  -0.80%  "
Set-TextValue $ws.Range("D40") "1.075"
$ws.Range("E40").Value = "  -0.33%  "
Set-TextValue $ws.Range("D41") "9.231"
$ws.Range("E41").Value = "  +2.43%  "
Set-TextValue $ws.Range("D42") "6.630"
$ws.Range("E42").Value = "  +1.13%  "
Set-TextValue $ws.Range("D43") "114.98"
$ws.Range("E43").Value = "  -1.07%  "
Set-TextValue $ws.Range("D44") "0.5084"
$ws.Range("E44").Value = "  +4.45%  "
$ws.Range("E45").Value = "  +0.68%  "
Set-TextValue $ws.Range("D47") "10.16"
$ws.Range("E47").Value = "  +0.77%  "
$ws.Range("E48").Value = "  +0.52%  "
Set-TextValue $ws.Range("D49") "37.91"
$ws.Range("E49").Value = "  -1.33%  "
Set-TextValue $ws.Range("D50") "0.06079"
$ws.Range("E50").Value = "  +2.52%  "
Set-TextValue $ws.Range("D51") "62.77"
$ws.Range("E51").Value = "  -1.46%  "
